$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we are about to (re)write keep plain-text storage,
# matching the original inline-string cells (avoids Excel auto-converting
# numeric-looking strings like "318.20" or "1.000" into real numbers).
$cellsToFormat = @("D2","D3","E3","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","E11","D12","D13","E13","D14","E14","D15","E15","D16","E16","E17","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","E27","D28","E28","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","B37","C37","D37","E37","B38","C38","D38","E38","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","E48","D49","E49","E50","D51","E51")
foreach ($addr in $cellsToFormat) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.722.07"
$ws.Range("D3").Value = "1.915.85"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").Value = "318.20"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").Value = "0.5197"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "0.3972"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "0.08522"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "42.72"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "6.301"
$ws.Range("D13").Value = "1.910.67"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").Value = "20.95"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "7.357"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "0.06754"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "17.98"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "6.040"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").Value = "29.705.10"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "2.212"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "2.127.27"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "159.37"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  +4.09%  "
$ws.Range("D30").Value = "128.45"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "1.088"
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("D32").Value = "0.1059"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "6.213"
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("D34").Value = "3.681"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "0.02500"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").Value = "0.06647"
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "9.150"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.250"
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("D39").Value = "0.2209"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "5.210"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("D41").Value = "0.6546"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "11.40"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "0.6147"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "13.32"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").Value = "3.691"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "125.01"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "78.55"
$ws.Range("E51").Value = "  +1.82%  "
